$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scholls = "Scholls  Position Catagories Add Functionality"
$schools = "Schools  Position Catagories Edit Functionality"
$passed  = "PASSED"
$failed  = "FAILED"
$chrome  = "chrome"

$rows = @(
    @($scholls, $passed, $chrome),
    @($schools, $failed, $chrome),
    @($schools, $failed, $chrome),
    @($scholls, $failed, $chrome),
    @($schools, $failed, $chrome),
    @($schools, $failed, $chrome),
    @($scholls, $passed, $chrome),
    @($schools, $passed, $chrome),
    @($schools, $failed, $chrome)
)

$startRow = 71
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
